$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text (European "."-grouped numbers are already
# inline strings in the source, e.g. "35.170.63"). Whenever the new price text would
# otherwise parse as a plain number (e.g. "40.38"), force the cell to Text format first
# so Excel does not silently convert it to a numeric value, then restore the default
# "Normal" style so no formatting residue is left on the cell.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.344"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0716"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0990"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "239.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0561"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.896"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0648"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0207"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "89.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.43"
$ws.Range("D51").Style = "Normal"

# Remaining coin / link / price / volume text updates. These are either pure text,
# multi-dot price strings, or space-padded percentages, none of which Excel
# reinterprets as numbers, so a direct .Value assignment is sufficient.
$ws.Range("D2").Value = "35.170.63"
$ws.Range("D3").Value = "1.893.13"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +7.98%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +10.98%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "2.170.36"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "1.891.00"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "35.223.65"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  -3.55%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +9.55%  "
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  -5.66%  "
$ws.Range("E38").Value = "  +10.98%  "
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E40").Value = "  +9.35%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E43").Value = "  +5.10%  "
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").Value = "1.336.63"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -8.35%  "
$ws.Range("E50").Value = "  -6.82%  "
$ws.Range("E51").Value = "  -3.63%  "
